$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 222.375
$ws.Range("I5").Value = 111.28571
$ws.Range("K5").Value = 111.28571
$ws.Range("M5").Value = 3.714290000000005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 78.454544
$ws.Range("I11").Value = 78.454544
$ws.Range("K11").Value = 78.454544
$ws.Range("M11").Value = 61.545456

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 362.66666
$ws.Range("I41").Value = 777
$ws.Range("J41").Value = 155.5
$ws.Range("K41").Value = 777
$ws.Range("L41").Value = 155.5
$ws.Range("M41").Value = -337
$ws.Range("N41").Value = -1035.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 18499.666
$ws.Range("I47").Value = 25250
$ws.Range("J47").Value = 4999
$ws.Range("K47").Value = 25250
$ws.Range("L47").Value = 4999
$ws.Range("M47").Value = -24278
$ws.Range("N47").Value = -6943

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6819.5
$ws.Range("J76").Value = 7099.4443
$ws.Range("L76").Value = 7099.4443
$ws.Range("N76").Value = -7729.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6819.5
$ws.Range("J79").Value = 7099.4443
$ws.Range("L79").Value = 7099.4443
$ws.Range("N79").Value = -9283.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4950.6
$ws.Range("I88").Value = 2375
$ws.Range("J88").Value = 5594.5
$ws.Range("K88").Value = 2375
$ws.Range("L88").Value = 5594.5
$ws.Range("M88").Value = -1969
$ws.Range("N88").Value = -6406.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4950.6
$ws.Range("I91").Value = 2375
$ws.Range("J91").Value = 5594.5
$ws.Range("K91").Value = 2375
$ws.Range("L91").Value = 5594.5
$ws.Range("M91").Value = -971
$ws.Range("N91").Value = -8402.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 1000000000
$ws.Range("I117").Value = 1000000000
$ws.Range("K117").Value = 1000000000
$ws.Range("M117").Value = -999995411

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2520.7334
$ws.Range("I132").Value = 2447
$ws.Range("K132").Value = 7341
$ws.Range("M132").Value = -4811

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3529.8667
$ws.Range("I137").Value = 1684.7142
$ws.Range("J137").Value = 5144.375
$ws.Range("K137").Value = 5054.142599999999
$ws.Range("L137").Value = 15433.125
$ws.Range("M137").Value = -2504.142599999999
$ws.Range("N137").Value = -20533.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1580.1428
$ws.Range("I45").Value = 1542.2
$ws.Range("K45").Value = 1542.2
$ws.Range("M45").Value = -1165.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4423.25
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4669.875
$ws.Range("J63").Value = 6782.375
$ws.Range("L63").Value = 6782.375
$ws.Range("N63").Value = -8154.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4669.875
$ws.Range("J66").Value = 6782.375
$ws.Range("L66").Value = 33911.875
$ws.Range("N66").Value = -40775.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4423.25
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -885

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4453.1055
$ws.Range("I105").Value = 3033.889
$ws.Range("K105").Value = 3033.889
$ws.Range("M105").Value = -1286.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 924999.5
$ws.Range("I122").Value = 924999
$ws.Range("K122").Value = 924999
$ws.Range("M122").Value = -920099

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200.78947
$ws.Range("J7").Value = 499.5
$ws.Range("L7").Value = 499.5
$ws.Range("N7").Value = -725.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 521.875
$ws.Range("I19").Value = 167.85715
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 167.85715
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 2.14285000000001
$ws.Range("N19").Value = -3340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 549.75
$ws.Range("I22").Value = 299.5
$ws.Range("K22").Value = 299.5
$ws.Range("M22").Value = 50.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 521.875
$ws.Range("I24").Value = 167.85715
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 167.85715
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = 2.14285000000001
$ws.Range("N24").Value = -3340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2820.8
$ws.Range("I122").Value = 2680.5625
$ws.Range("K122").Value = 8041.6875
$ws.Range("M122").Value = -5591.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2204.95
$ws.Range("I132").Value = 2081.625
$ws.Range("K132").Value = 6244.875
$ws.Range("M132").Value = -3714.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50048.25
$ws.Range("J2").Value = 91.28571
$ws.Range("L2").Value = 547.71426
$ws.Range("N2").Value = -773.71426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2002.2307
$ws.Range("J34").Value = 2488.6667
$ws.Range("L34").Value = 7466.000100000001
$ws.Range("N34").Value = -7634.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1344.5714
$ws.Range("J63").Value = 2250
$ws.Range("L63").Value = 6750
$ws.Range("N63").Value = -8248

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 1344.5714
$ws.Range("J66").Value = 2250
$ws.Range("L66").Value = 20250
$ws.Range("N66").Value = -27738

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1027.7
$ws.Range("J113").Value = 1283.5714
$ws.Range("L113").Value = 3850.7142
$ws.Range("N113").Value = -8190.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 2624.8
$ws.Range("I119").Value = 2624.8
$ws.Range("K119").Value = 7874.400000000001
$ws.Range("M119").Value = -3036.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 3979888.8
$ws.Range("I128").Value = 3979888.8
$ws.Range("K128").Value = 11939666.4
$ws.Range("M128").Value = -11934686.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 10000
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 261.8889
$ws.Range("I2").Value = 60.727272
$ws.Range("J2").Value = 578
$ws.Range("K2").Value = 60.727272
$ws.Range("L2").Value = 578
$ws.Range("M2").Value = 52.272728
$ws.Range("N2").Value = -804

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3933
$ws.Range("I70").Value = 3933
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3933
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3663
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 3933
$ws.Range("I73").Value = 3933
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3933
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2997
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4372.067
$ws.Range("J113").Value = 4157.9165
$ws.Range("L113").Value = 4157.9165
$ws.Range("N113").Value = -8497.9165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14374.75
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 14374.75
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2917.7144
$ws.Range("I46").Value = 1471.5
$ws.Range("K46").Value = 1471.5
$ws.Range("M46").Value = -1283.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 201
$ws.Range("I55").Value = 189.375
$ws.Range("J55").Value = 294
$ws.Range("K55").Value = 189.375
$ws.Range("L55").Value = 294
$ws.Range("M55").Value = -16.375
$ws.Range("N55").Value = -640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3003.92
$ws.Range("I82").Value = 3263.2222
$ws.Range("J82").Value = 2337.1428
$ws.Range("K82").Value = 3263.2222
$ws.Range("L82").Value = 2337.1428
$ws.Range("M82").Value = -2902.2222
$ws.Range("N82").Value = -3059.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3003.92
$ws.Range("I85").Value = 3263.2222
$ws.Range("J85").Value = 2337.1428
$ws.Range("K85").Value = 3263.2222
$ws.Range("L85").Value = 2337.1428
$ws.Range("M85").Value = -2015.2222
$ws.Range("N85").Value = -4833.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2142.8572
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2166.6667
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3248.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4733.25
$ws.Range("I132").Value = 1799.75
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 5399.25
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -2869.25
$ws.Range("N132").Value = -23660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
